$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.220024228096008
$ws.Range("B1").Value = 2.699231624603271
$ws.Range("C1").Value = 4.47798490524292
$ws.Range("D1").Value = 2.138338327407837
$ws.Range("E1").Value = 1.161771059036255
